# Commit: "added checks for no description, invalid header"
#
# The "params" sheet gains a new "description" column (column M, after the
# existing "ui variable" column L) so that downstream validation code can
# check rows for a missing description / an invalid header. Only the
# header row needs the new label - rows 2 and 3 have no value for it yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

# New header cell - becomes a new shared-string entry ("description") and
# extends the sheet's used range/dimension out to column M.
$ws.Range("M1").Value = "description"

# Move/resize the selection onto the newly added header cell, matching the
# updated <selection activeCell="M1" sqref="M1"/> in the sheet view.
$ws.Range("M1").Select()
